# EmailNoMatch.xlsx - RAD test data update
# Updates the "Execution Date" timestamps (column B, rows 2-33) on Sheet1
# with the results of the latest Katalon test run
# ("FEIN/SSN Object Identification").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "Wed Nov 01 15:45:59 EDT 2023"
$ws.Range("B3").Value  = "Wed Nov 01 15:46:11 EDT 2023"
$ws.Range("B4").Value  = "Wed Nov 01 15:46:22 EDT 2023"
$ws.Range("B5").Value  = "Wed Nov 01 15:46:34 EDT 2023"
$ws.Range("B6").Value  = "Wed Nov 01 15:46:46 EDT 2023"
$ws.Range("B7").Value  = "Wed Nov 01 15:46:58 EDT 2023"
$ws.Range("B8").Value  = "Wed Nov 01 15:47:10 EDT 2023"
$ws.Range("B9").Value  = "Wed Nov 01 15:47:22 EDT 2023"
$ws.Range("B10").Value = "Wed Nov 01 15:47:34 EDT 2023"
$ws.Range("B11").Value = "Wed Nov 01 15:47:46 EDT 2023"
$ws.Range("B12").Value = "Wed Nov 01 15:47:58 EDT 2023"
$ws.Range("B13").Value = "Wed Nov 01 15:48:10 EDT 2023"
$ws.Range("B14").Value = "Wed Nov 01 15:48:22 EDT 2023"
$ws.Range("B15").Value = "Wed Nov 01 15:48:34 EDT 2023"
$ws.Range("B16").Value = "Wed Nov 01 15:48:46 EDT 2023"
$ws.Range("B17").Value = "Wed Nov 01 15:48:58 EDT 2023"
$ws.Range("B18").Value = "Wed Nov 01 15:49:10 EDT 2023"
$ws.Range("B19").Value = "Wed Nov 01 15:49:22 EDT 2023"
$ws.Range("B20").Value = "Wed Nov 01 15:49:34 EDT 2023"
$ws.Range("B21").Value = "Wed Nov 01 15:49:45 EDT 2023"
$ws.Range("B22").Value = "Wed Nov 01 15:49:57 EDT 2023"
$ws.Range("B23").Value = "Wed Nov 01 15:50:09 EDT 2023"
$ws.Range("B24").Value = "Wed Nov 01 15:50:21 EDT 2023"
$ws.Range("B25").Value = "Wed Nov 01 15:50:33 EDT 2023"
$ws.Range("B26").Value = "Wed Nov 01 15:50:45 EDT 2023"
$ws.Range("B27").Value = "Wed Nov 01 15:50:57 EDT 2023"
$ws.Range("B28").Value = "Wed Nov 01 15:51:09 EDT 2023"
$ws.Range("B29").Value = "Wed Nov 01 15:51:22 EDT 2023"
$ws.Range("B30").Value = "Wed Nov 01 15:51:34 EDT 2023"
$ws.Range("B31").Value = "Wed Nov 01 15:51:46 EDT 2023"
$ws.Range("B32").Value = "Wed Nov 01 15:51:58 EDT 2023"
$ws.Range("B33").Value = "Wed Nov 01 15:52:10 EDT 2023"
